$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Done" (column D) x-marks, new column E header, and "Comment" (column E)
# notes for the first block of tasks (rows 2-8), written in the order that
# reproduces the workbook's shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "x"
$ws.Range("E2").Value = "Comment"
$ws.Range("E6").Value = "Setorder market vaue to June"
$ws.Range("E4").Value = "ongoig"
$ws.Range("D5").Value = "x"
$ws.Range("E7").Value = "ongoig"
$ws.Range("D8").Value = "x"

# ---------------------------------------------------------------------------
# Three brand-new task rows get appended at the bottom of the list first,
# before the existing rows 9-11 are edited in place (this reproduces the
# shared-string write order captured in the workbook).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Task 10"
$ws.Range("A13").Value = "Task 11"
$ws.Range("A14").Value = "Task 12"
$ws.Range("A15").Value = "Task 13"
$ws.Range("A16").Value = "Task 14"
$ws.Range("A17").Value = "Task 15"

# ---------------------------------------------------------------------------
# Row 10 becomes the new "Strategy" note (wrapped, taller row), row 9 gets
# a new description + comment.
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "Strategy: Mean Variance, min Vol, max sharpe, … to be informed and be able to aks questions on 20 May"
$ws.Range("B10").Value = "All"
$ws.Range("C9").Value = "Begin in 1994"
$ws.Range("E9").Value = "From Hanauer Mail"

# The old descriptions of rows 9-11 (Task 7/8/9) slide down into the newly
# appended rows, and row 11 picks up the old row-9 description.
$ws.Range("C14").Value = "Cross-section regression"
$ws.Range("C12").Value = "Run CAPM regression on hedge portfolio to check for excess return"
$ws.Range("C13").Value = "Re-Create Table 3 and Table 4 from Hanauer, Lauterbach Paper"
$ws.Range("C11").Value = "Create hedge portfolio from factors"

# Wrap the long strategy text and make the row tall enough to show it.
$ws.Range("C10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 32

# ---------------------------------------------------------------------------
# Column E width so the comments are readable
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 36.33

# ---------------------------------------------------------------------------
# Match the saved selection / active cell
# ---------------------------------------------------------------------------
$null = $ws.Range("E8").Select()
